# Daily attendance processing - reorder "Recorded By" names so that
# "Administrator" moves from the front of the list to the end, e.g.:
#   "Administrator, Miss Dina Nasr, Developer" -> "Miss Dina Nasr, Developer, Administrator"
#   "Administrator, Miss Dina Nasr"            -> "Miss Dina Nasr, Administrator"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "Administrator, *") {
        $rest = $val.Substring("Administrator, ".Length)
        $cell.Value = "$rest, Administrator"
    }
}
